# Candidate workbook update:
#  - "position" column split into 3 columns: position_applied_for_1/2/3
#  - new position values populated per row
#  - row 5 email changed to das@yopmail.com
#  - related column width / selection metadata updates

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert two new blank columns at W (23), pushing old W/X/Y to Y/Z/AA ---
$ws.Columns.Item(23).Insert()
$ws.Columns.Item(23).Insert()

# --- Header row ---
$ws.Range("V1").Value = "position_applied_for_1"
$ws.Range("W1").Value = "position_applied_for_2"
$ws.Range("X1").Value = "position_applied_for_3"

# --- Column V: position_applied_for_1 (filled top-to-bottom) ---
$ws.Range("V2").Value = "ACCOUNTANT"
$ws.Range("V3").Value = "ACCOUNTANT TALLY"
$ws.Range("V4").Value = "ALUMINIUM FABRICATOR"
$ws.Range("V5").Value = "ANIMAL WARDEN"

# --- Column W: position_applied_for_2 (filled top-to-bottom) ---
$ws.Range("W2").Value = "ANY HELPER"
$ws.Range("W3").Value = "ARBIC CHEF"
$ws.Range("W4").Value = "AREA RESTURANT MANAGER"
$ws.Range("W5").Value = "ARGON WELDER"

# --- Column X: position_applied_for_3 (filled top-to-bottom) ---
$ws.Range("X2").Value = "ASST. COOK"
$ws.Range("X3").Value = "ASST. COOK TANDOOR"
$ws.Range("X4").Value = "ASST. ELECTRICIAN"
$ws.Range("X5").Value = "ASST. INDIAN COOK"

# --- Row 5: updated candidate email ---
$ws.Range("L5").Value = "das@yopmail.com"

# --- Column widths for the re-shuffled / new columns ---
$ws.Columns.Item(22).ColumnWidth = 21.17   # V  -> width 22
$ws.Columns.Item(23).ColumnWidth = 21.17   # W  -> width 22
$ws.Columns.Item(24).ColumnWidth = 21.17   # X  -> width 22
$ws.Columns.Item(25).ColumnWidth = 15.83   # Y  -> width ~16.71 (closest reachable: 16.667)
$ws.Columns.Item(26).ColumnWidth = 15.83   # Z  -> width ~16.71 (closest reachable: 16.667)
$ws.Columns.Item(27).ColumnWidth = 10.6    # AA -> width ~11.43 (closest reachable: 11.5)

# --- View / selection state ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 2    # topLeftCell = B1
$win.ScrollRow = 1
$ws.Range("H16").Select()
